$wb = $excel.ActiveWorkbook

# --- farm_layout sheet: insert a new "cocoa / b" planting row, shifting the
#     existing "coconut" row down from row 3 to row 4. ----------------------
$wsFarm = $wb.Worksheets.Item("farm_layout")

# Insert a new blank row at row 3 (pushes old row 3 "coconut" down to row 4).
$wsFarm.Rows.Item(3).Insert()

# Populate the freshly inserted row 3 with the new planting entry.
$wsFarm.Range("A3").Value = "cocoa"
$wsFarm.Range("B3").Value = "b"
$wsFarm.Range("C3").Value = 1
$wsFarm.Range("D3").Value = 5
$wsFarm.Range("E3").Value = 0
$wsFarm.Range("F3").Value = 0
$wsFarm.Range("G3").Value = 5000

# --- simulation_control sheet keeps its own selection (A7), it just stops
#     being the active tab once another sheet becomes active below. --------
$wsSim = $wb.Worksheets.Item("simulation_control")
$wsSim.Range("A7").Select()

# --- farm_layout becomes the active/selected sheet, with I18 selected. ----
$wsFarm.Range("I18").Select()
